# Convert risk difference from decimal to percent in benefit tables
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-CellText($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # Shrink the range so it doesn't include the trailing end-of-cell /
    # paragraph-mark characters.
    $rng.End = $rng.End - 1
    # wdReplaceOne (1), not wdReplaceAll, and Wrap=wdFindStop (0) so the
    # search/replace stays confined to this single cell.
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                $true, 0, $false, $newText, 1)
    if (-not $found) {
        Write-Host "NOT FOUND: row=$row col=$col old='$oldText'"
    }
}

# Header: "Posterior Median RD (95% CrI)" -> "Posterior Median RD, % (95% CrI)"
Replace-CellText 2 3 " RD (95% " " RD, % (95% "

# Optimistic / Weak
Replace-CellText 4 3 "0.05 (-0.02, 0.12)" "5.2 (-1.6, 12.1)"

# Optimistic / Moderate
Replace-CellText 5 3 "0.04 (-0.01, 0.09)" "4.0 (-0.7, 8.8)"
Replace-CellText 5 5 "79" "80"

# Optimistic / Strong
Replace-CellText 6 3 "0.03 (0, 0.07)" "3.4 (0.2, 6.8)"
Replace-CellText 6 5 "81" "80"

# Neutral / Weak
Replace-CellText 7 3 "0.07 (-0.02, 0.15)" "6.6 (-2.1, 15.3)"
Replace-CellText 7 4 "94" "93"
Replace-CellText 7 7 "23" "22"

# Neutral / Moderate
Replace-CellText 8 3 "0.03 (-0.03, 0.10)" "3.5 (-3.0, 9.8)"

# Neutral / Strong
Replace-CellText 9 3 "0.02 (-0.03, 0.06)" "1.8 (-2.8, 6.3)"
Replace-CellText 9 5 "45" "46"

# Pessimistic / Weak
Replace-CellText 10 3 "0.03 (-0.04, 0.10)" "2.9 (-4.0, 9.8)"
Replace-CellText 10 6 "28" "27"

# Pessimistic / Moderate
Replace-CellText 11 3 "0 (-0.05, 0.05)" "-0.2 (-5.0, 4.5)"

# Pessimistic / Strong
Replace-CellText 12 3 "-0.02 (-0.05, 0.02)" "-1.6 (-5.0, 1.6)"
Replace-CellText 12 5 "1" "2"
